$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.170.31'
$ws.Range("E2").Value = '  -2.18%  '
$ws.Range("D3").Value = '1.577.04'
$ws.Range("E3").Value = '  -1.68%  '
$ws.Range("E4").Value = '  -0.50%  '
$ws.Range("D5").Value = '''208.79'
$ws.Range("E5").Value = '  -1.46%  '
$ws.Range("E6").Value = '  -3.07%  '
$ws.Range("E7").Value = '  -0.48%  '
$ws.Range("E8").Value = '  -1.59%  '
$ws.Range("D9").Value = '''0.245'
$ws.Range("E9").Value = '  -0.94%  '
$ws.Range("E10").Value = '  -0.39%  '
$ws.Range("D11").Value = '''0.0843'
$ws.Range("E11").Value = '  -0.53%  '
$ws.Range("D12").Value = '1.798.60'
$ws.Range("E12").Value = '  -1.66%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.591.94'
$ws.Range("E13").Value = '  -0.74%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '''4.05'
$ws.Range("E14").Value = '  -0.49%  '
$ws.Range("D15").Value = '''0.515'
$ws.Range("E15").Value = '  -2.11%  '
$ws.Range("D16").Value = '''64.39'
$ws.Range("E16").Value = '  -1.10%  '
$ws.Range("D17").Value = '26.159.33'
$ws.Range("E17").Value = '  -2.16%  '
$ws.Range("D18").Value = '0.0₃0727'
$ws.Range("E18").Value = '  -2.11%  '
$ws.Range("D19").Value = '''7.31'
$ws.Range("E19").Value = '  +1.37%  '
$ws.Range("D20").Value = '''208.62'
$ws.Range("E20").Value = '  -0.52%  '
$ws.Range("E22").Value = '  -1.31%  '
$ws.Range("D23").Value = '''2.16'
$ws.Range("E23").Value = '  -2.80%  '
$ws.Range("E24").Value = '  -2.36%  '
$ws.Range("D25").Value = '''144.20'
$ws.Range("E25").Value = '  +0.29%  '
$ws.Range("E26").Value = '  -0.35%  '
$ws.Range("E27").Value = '  -1.93%  '
$ws.Range("E28").Value = '  -1.80%  '
$ws.Range("D29").Value = '''15.23'
$ws.Range("E29").Value = '  -0.90%  '
$ws.Range("E30").Value = '  -0.13%  '
$ws.Range("E31").Value = '  -1.27%  '
$ws.Range("E33").Value = '  +0.78%  '
$ws.Range("D34").Value = '1.280.79'
$ws.Range("E34").Value = '  -0.52%  '
$ws.Range("E35").Value = '  +4.03%  '
$ws.Range("E36").Value = '  -1.68%  '
$ws.Range("E37").Value = '  -1.37%  '
$ws.Range("D38").Value = '''0.0166'
$ws.Range("E38").Value = '  -2.42%  '
$ws.Range("D39").Value = '''1.10'
$ws.Range("E39").Value = '  -10.37%  '
$ws.Range("E40").Value = '  -2.17%  '
$ws.Range("E41").Value = '  -0.49%  '
$ws.Range("D42").Value = '''5.58'
$ws.Range("E42").Value = '  +2.59%  '
$ws.Range("E43").Value = '  -1.70%  '
$ws.Range("E44").Value = '  -3.18%  '
$ws.Range("D45").Value = '''62.28'
$ws.Range("E45").Value = '  -0.69%  '
$ws.Range("D46").Value = '1.711.74'
$ws.Range("E46").Value = '  -1.64%  '
$ws.Range("D47").Value = '''88.63'
$ws.Range("E47").Value = '  -1.97%  '
$ws.Range("E48").Value = '  -2.57%  '
$ws.Range("E49").Value = '  -4.26%  '
$ws.Range("D50").Value = '''0.100'
$ws.Range("E50").Value = '  -1.85%  '
$ws.Range("E51").Value = '  -1.55%  '

# Reset style to Normal for text-forced numeric-looking cells so the
# quote-prefix indicator does not introduce a spurious style (matches source).
$ws.Range("D5").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").Style = "Normal"
